# Auto-generated edit script: updates profit-calculation cells (columns H-N)
# across multiple Leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to
# reflect refreshed market-price data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 1635.5834  # H80
$ws.Cells.Item(80, 9).Value = 1763  # I80
$ws.Cells.Item(80, 11).Value = 5289  # K80
$ws.Cells.Item(80, 13).Value = -4291  # M80
$ws.Cells.Item(83, 8).Value = 1635.5834  # H83
$ws.Cells.Item(83, 9).Value = 1763  # I83
$ws.Cells.Item(83, 11).Value = 15867  # K83
$ws.Cells.Item(83, 13).Value = -10875  # M83
$ws.Cells.Item(92, 8).Value = 974.9643  # H92
$ws.Cells.Item(92, 9).Value = 785.8889  # I92
$ws.Cells.Item(92, 10).Value = 1315.3  # J92
$ws.Cells.Item(92, 11).Value = 785.8889  # K92
$ws.Cells.Item(92, 12).Value = 1315.3  # L92
$ws.Cells.Item(92, 13).Value = 462.1111  # M92
$ws.Cells.Item(92, 14).Value = -3811.3  # N92
$ws.Cells.Item(105, 8).Value = 89999  # H105
$ws.Cells.Item(105, 10).Value = 0  # J105
$ws.Cells.Item(105, 12).Value = 0  # L105
$ws.Cells.Item(105, 14).ClearContents() | Out-Null  # N105
$ws.Cells.Item(113, 8).Value = 5349.7144  # H113
$ws.Cells.Item(113, 9).Value = 4799.5  # I113
$ws.Cells.Item(113, 11).Value = 4799.5  # K113
$ws.Cells.Item(113, 13).Value = -1545.5  # M113
$ws.Cells.Item(132, 8).Value = 4824.894  # H132
$ws.Cells.Item(132, 9).Value = 2586.8794  # I132
$ws.Cells.Item(132, 11).Value = 7760.638199999999  # K132
$ws.Cells.Item(132, 13).Value = -5230.638199999999  # M132
$ws.Cells.Item(135, 8).Value = 2111  # H135
$ws.Cells.Item(135, 9).Value = 1617.9375  # I135
$ws.Cells.Item(135, 11).Value = 14561.4375  # K135
$ws.Cells.Item(135, 13).Value = -12026.4375  # M135
$ws.Cells.Item(137, 8).Value = 4843.9346  # H137
$ws.Cells.Item(137, 9).Value = 5393.0293  # I137
$ws.Cells.Item(137, 10).Value = 3288.1667  # J137
$ws.Cells.Item(137, 11).Value = 16179.0879  # K137
$ws.Cells.Item(137, 12).Value = 9864.500100000001  # L137
$ws.Cells.Item(137, 13).Value = -13629.0879  # M137
$ws.Cells.Item(137, 14).Value = -14964.5001  # N137
$ws.Cells.Item(138, 8).Value = 2527.5557  # H138
$ws.Cells.Item(138, 9).Value = 1490.4762  # I138
$ws.Cells.Item(138, 11).Value = 4471.4286  # K138
$ws.Cells.Item(138, 13).Value = 668.5713999999998  # M138
$ws.Cells.Item(141, 8).Value = 15362.75  # H141
$ws.Cells.Item(141, 9).Value = 15362.75  # I141
$ws.Cells.Item(141, 11).Value = 46088.25  # K141
$ws.Cells.Item(141, 13).Value = -40908.25  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7529.1455  # H32
$ws.Cells.Item(32, 9).Value = 6039.698  # I32
$ws.Cells.Item(32, 10).Value = 46999.5  # J32
$ws.Cells.Item(32, 11).Value = 6039.698  # K32
$ws.Cells.Item(32, 12).Value = 46999.5  # L32
$ws.Cells.Item(32, 13).Value = -5752.698  # M32
$ws.Cells.Item(32, 14).Value = -47573.5  # N32
$ws.Cells.Item(132, 8).Value = 435.65  # H132
$ws.Cells.Item(132, 9).Value = 392.78946  # I132
$ws.Cells.Item(132, 10).Value = 1250  # J132
$ws.Cells.Item(132, 11).Value = 1178.36838  # K132
$ws.Cells.Item(132, 12).Value = 3750  # L132
$ws.Cells.Item(132, 13).Value = 1351.63162  # M132
$ws.Cells.Item(132, 14).Value = -8810  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(93, 8).Value = 0  # H93
$ws.Cells.Item(93, 10).Value = 0  # J93
$ws.Cells.Item(93, 12).Value = 0  # L93
$ws.Cells.Item(93, 14).ClearContents() | Out-Null  # N93

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(87, 8).Value = 0  # H87
$ws.Cells.Item(87, 9).Value = 0  # I87
$ws.Cells.Item(87, 11).Value = 0  # K87
$ws.Cells.Item(87, 13).ClearContents() | Out-Null  # M87
$ws.Cells.Item(90, 8).Value = 0  # H90
$ws.Cells.Item(90, 9).Value = 0  # I90
$ws.Cells.Item(90, 11).Value = 0  # K90
$ws.Cells.Item(90, 13).ClearContents() | Out-Null  # M90
$ws.Cells.Item(105, 8).Value = 1765.5  # H105
$ws.Cells.Item(105, 9).Value = 943.6667  # I105
$ws.Cells.Item(105, 10).Value = 2998.25  # J105
$ws.Cells.Item(105, 11).Value = 943.6667  # K105
$ws.Cells.Item(105, 12).Value = 2998.25  # L105
$ws.Cells.Item(105, 13).Value = 803.3333  # M105
$ws.Cells.Item(105, 14).Value = -6492.25  # N105

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 792.6667  # H5
$ws.Cells.Item(5, 9).Value = 792.6667  # I5
$ws.Cells.Item(5, 10).Value = 0  # J5
$ws.Cells.Item(5, 11).Value = 2378.0001  # K5
$ws.Cells.Item(5, 12).Value = 0  # L5
$ws.Cells.Item(5, 13).Value = -2266.0001  # M5
$ws.Cells.Item(5, 14).ClearContents() | Out-Null  # N5
$ws.Cells.Item(9, 8).Value = 84169340  # H9
$ws.Cells.Item(9, 10).Value = 100003110  # J9
$ws.Cells.Item(9, 12).Value = 300009330  # L9
$ws.Cells.Item(9, 14).Value = -300009778  # N9
$ws.Cells.Item(12, 8).Value = 121  # H12
$ws.Cells.Item(12, 9).Value = 79.40000000000001  # I12
$ws.Cells.Item(12, 10).Value = 134.86667  # J12
$ws.Cells.Item(12, 11).Value = 238.2  # K12
$ws.Cells.Item(12, 12).Value = 404.60001  # L12
$ws.Cells.Item(12, 13).Value = -65.20000000000002  # M12
$ws.Cells.Item(12, 14).Value = -750.60001  # N12
$ws.Cells.Item(34, 8).Value = 913.55554  # H34
$ws.Cells.Item(34, 9).Value = 100  # I34
$ws.Cells.Item(34, 10).Value = 1146  # J34
$ws.Cells.Item(34, 11).Value = 300  # K34
$ws.Cells.Item(34, 12).Value = 3438  # L34
$ws.Cells.Item(34, 13).Value = -216  # M34
$ws.Cells.Item(34, 14).Value = -3606  # N34
$ws.Cells.Item(39, 8).Value = 12075.77  # H39
$ws.Cells.Item(39, 9).Value = 1499.5  # I39
$ws.Cells.Item(39, 10).Value = 13998.728  # J39
$ws.Cells.Item(39, 11).Value = 4498.5  # K39
$ws.Cells.Item(39, 12).Value = 41996.18399999999  # L39
$ws.Cells.Item(39, 13).Value = -4204.5  # M39
$ws.Cells.Item(39, 14).Value = -42584.18399999999  # N39
$ws.Cells.Item(55, 8).Value = 2987.0386  # H55
$ws.Cells.Item(55, 9).Value = 470.93332  # I55
$ws.Cells.Item(55, 10).Value = 6418.091  # J55
$ws.Cells.Item(55, 11).Value = 1412.79996  # K55
$ws.Cells.Item(55, 12).Value = 19254.273  # L55
$ws.Cells.Item(55, 13).Value = -1235.79996  # M55
$ws.Cells.Item(55, 14).Value = -19608.273  # N55
$ws.Cells.Item(108, 8).Value = 1050.5  # H108
$ws.Cells.Item(108, 9).Value = 1050.5  # I108
$ws.Cells.Item(108, 11).Value = 3151.5  # K108
$ws.Cells.Item(108, 13).Value = -271.5  # M108
$ws.Cells.Item(135, 8).Value = 792.6667  # H135
$ws.Cells.Item(135, 9).Value = 792.6667  # I135
$ws.Cells.Item(135, 10).Value = 0  # J135
$ws.Cells.Item(135, 11).Value = 7134.0003  # K135
$ws.Cells.Item(135, 12).Value = 0  # L135
$ws.Cells.Item(135, 13).Value = -4599.0003  # M135
$ws.Cells.Item(135, 14).ClearContents() | Out-Null  # N135
$ws.Cells.Item(140, 8).Value = 844.8889  # H140
$ws.Cells.Item(140, 9).Value = 901.75  # I140
$ws.Cells.Item(140, 11).Value = 2705.25  # K140
$ws.Cells.Item(140, 13).Value = 2474.75  # M140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2837.6785  # H122
$ws.Cells.Item(122, 9).Value = 2294.9546  # I122
$ws.Cells.Item(122, 10).Value = 4827.6665  # J122
$ws.Cells.Item(122, 11).Value = 6884.8638  # K122
$ws.Cells.Item(122, 12).Value = 14482.9995  # L122
$ws.Cells.Item(122, 13).Value = -4434.8638  # M122
$ws.Cells.Item(122, 14).Value = -19382.9995  # N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3406.6428  # H46
$ws.Cells.Item(46, 9).Value = 1075  # I46
$ws.Cells.Item(46, 11).Value = 1075  # K46
$ws.Cells.Item(46, 13).Value = -887  # M46
$ws.Cells.Item(55, 8).Value = 522.375  # H55
$ws.Cells.Item(55, 9).Value = 512  # I55
$ws.Cells.Item(55, 10).Value = 567.3333  # J55
$ws.Cells.Item(55, 11).Value = 512  # K55
$ws.Cells.Item(55, 12).Value = 567.3333  # L55
$ws.Cells.Item(55, 13).Value = -339  # M55
$ws.Cells.Item(55, 14).Value = -913.3333  # N55
$ws.Cells.Item(82, 8).Value = 2265.7937  # H82
$ws.Cells.Item(82, 9).Value = 2269.2842  # I82
$ws.Cells.Item(82, 11).Value = 2269.2842  # K82
$ws.Cells.Item(82, 13).Value = -1908.2842  # M82
$ws.Cells.Item(85, 8).Value = 2265.7937  # H85
$ws.Cells.Item(85, 9).Value = 2269.2842  # I85
$ws.Cells.Item(85, 11).Value = 2269.2842  # K85
$ws.Cells.Item(85, 13).Value = -1021.2842  # M85
$ws.Cells.Item(118, 8).Value = 76535.39999999999  # H118
$ws.Cells.Item(118, 10).Value = 76535.39999999999  # J118
$ws.Cells.Item(118, 12).Value = 76535.39999999999  # L118
$ws.Cells.Item(118, 14).Value = -79849.39999999999  # N118

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(116, 8).Value = 99949.5  # H116
$ws.Cells.Item(116, 10).Value = 99949.5  # J116
$ws.Cells.Item(116, 12).Value = 99949.5  # L116
$ws.Cells.Item(116, 14).Value = -109127.5  # N116
$ws.Cells.Item(132, 8).Value = 3937.389  # H132
$ws.Cells.Item(132, 9).Value = 4508.154  # I132
$ws.Cells.Item(132, 11).Value = 13524.462  # K132
$ws.Cells.Item(132, 13).Value = -10994.462  # M132
$ws.Cells.Item(136, 8).Value = 1252.4  # H136
$ws.Cells.Item(136, 9).Value = 1254.5834  # I136
$ws.Cells.Item(136, 11).Value = 3763.7502  # K136
$ws.Cells.Item(136, 13).Value = -1213.7502  # M136

